$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.515.98'
$ws.Cells.Item(2, 5).Value = '  +0.22%  '

$ws.Cells.Item(3, 4).Value = '1.789.85'
$ws.Cells.Item(3, 5).Value = '  -2.05%  '

$ws.Cells.Item(4, 5).Value = '  -0.17%  '

$ws.Cells.Item(5, 4).Value = '225.45'
$ws.Cells.Item(5, 5).Value = '  -2.16%  '

$ws.Cells.Item(6, 4).Value = '0.554'
$ws.Cells.Item(6, 5).Value = '  -3.68%  '

$ws.Cells.Item(7, 4).Value = '0.998'
$ws.Cells.Item(7, 5).Value = '  -0.19%  '

$ws.Cells.Item(8, 4).Value = '33.26'
$ws.Cells.Item(8, 5).Value = '  +5.38%  '

$ws.Cells.Item(9, 4).Value = '0.283'
$ws.Cells.Item(9, 5).Value = '  -1.92%  '

$ws.Cells.Item(10, 4).Value = '0.0664'
$ws.Cells.Item(10, 5).Value = '  -2.19%  '

$ws.Cells.Item(11, 4).Value = '0.0931'
$ws.Cells.Item(11, 5).Value = '  -0.21%  '

$ws.Cells.Item(12, 4).Value = '2.048.10'
$ws.Cells.Item(12, 5).Value = '  -2.01%  '

$ws.Cells.Item(13, 4).Value = '11.10'
$ws.Cells.Item(13, 5).Value = '  +7.43%  '

$ws.Cells.Item(14, 4).Value = '1.799.41'
$ws.Cells.Item(14, 5).Value = '  -1.61%  '

$ws.Cells.Item(15, 4).Value = '0.637'
$ws.Cells.Item(15, 5).Value = '  -2.94%  '

$ws.Cells.Item(16, 4).Value = '34.455.71'
$ws.Cells.Item(16, 5).Value = '  +0.19%  '

$ws.Cells.Item(17, 4).Value = '4.26'
$ws.Cells.Item(17, 5).Value = '  -1.38%  '

$ws.Cells.Item(18, 4).Value = '69.32'
$ws.Cells.Item(18, 5).Value = '  -1.74%  '

$ws.Cells.Item(19, 4).Value = '256.93'
$ws.Cells.Item(19, 5).Value = '  -0.62%  '

$ws.Cells.Item(20, 4).Value = '0.0₃0751'
$ws.Cells.Item(20, 5).Value = '  -1.16%  '

$ws.Cells.Item(21, 4).Value = '0.999'
$ws.Cells.Item(21, 5).Value = '  -0.13%  '

$ws.Cells.Item(22, 4).Value = '10.47'
$ws.Cells.Item(22, 5).Value = '  -1.78%  '

$ws.Cells.Item(23, 4).Value = '4.22'
$ws.Cells.Item(23, 5).Value = '  -2.99%  '

$ws.Cells.Item(24, 5).Value = '  -4.14%  '

$ws.Cells.Item(25, 4).Value = '157.70'
$ws.Cells.Item(25, 5).Value = '  -1.20%  '

$ws.Cells.Item(26, 4).Value = '16.52'
$ws.Cells.Item(26, 5).Value = '  -1.79%  '

$ws.Cells.Item(27, 4).Value = '7.07'
$ws.Cells.Item(27, 5).Value = '  -1.52%  '

$ws.Cells.Item(28, 5).Value = '  -3.69%  '

$ws.Cells.Item(29, 5).Value = '  -0.20%  '

$ws.Cells.Item(30, 4).Value = '3.80'
$ws.Cells.Item(30, 5).Value = '  -2.51%  '

$ws.Cells.Item(31, 5).Value = '  -2.09%  '

$ws.Cells.Item(32, 5).Value = '  -2.11%  '

$ws.Cells.Item(33, 4).Value = '3.61'
$ws.Cells.Item(33, 5).Value = '  +0.31%  '

$ws.Cells.Item(34, 5).Value = '  +5.50%  '

$ws.Cells.Item(35, 4).Value = '1.455.94'
$ws.Cells.Item(35, 5).Value = '  -5.78%  '

$ws.Cells.Item(36, 5).Value = '  -1.93%  '

$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(37, 4).Value = '0.632'
$ws.Cells.Item(37, 5).Value = '  -1.13%  '

$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(38, 4).Value = '0.0189'
$ws.Cells.Item(38, 5).Value = '  -1.36%  '

$ws.Cells.Item(39, 4).Value = '84.23'

$ws.Cells.Item(40, 5).Value = '  +1.41%  '

$ws.Cells.Item(41, 5).Value = '  -0.72%  '

$ws.Cells.Item(42, 4).Value = '0.896'
$ws.Cells.Item(42, 5).Value = '  -2.18%  '

$ws.Cells.Item(43, 4).Value = '2.09'
$ws.Cells.Item(43, 5).Value = '  -2.56%  '

$ws.Cells.Item(44, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(44, 4).Value = '12.72'
$ws.Cells.Item(44, 5).Value = '  +4.03%  '

$ws.Cells.Item(45, 2).Value = 'Kaspa'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(45, 4).Value = '0.0507'
$ws.Cells.Item(45, 5).Value = '  -4.12%  '

$ws.Cells.Item(46, 2).Value = 'FraxShare'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, 4).Value = '5.92'
$ws.Cells.Item(46, 5).Value = '  +1.63%  '

$ws.Cells.Item(47, 2).Value = 'WEMIXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(47, 4).Value = '1.05'
$ws.Cells.Item(47, 5).Value = '  -2.58%  '

$ws.Cells.Item(48, 4).Value = '1.947.55'
$ws.Cells.Item(48, 5).Value = '  -1.52%  '

$ws.Cells.Item(49, 4).Value = '0.998'
$ws.Cells.Item(49, 5).Value = '  -0.22%  '

$ws.Cells.Item(50, 4).Value = '99.34'
$ws.Cells.Item(50, 5).Value = '  -0.09%  '

$ws.Cells.Item(51, 4).Value = '50.61'
$ws.Cells.Item(51, 5).Value = '  -2.57%  '
